# Normalise the comma-separated "Recorded By" entries in column G so the
# automated "System"/"system" actor is no longer listed first — e.g.
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".
#
# Applied by reversing the order of the comma-separated list for every
# multi-entry cell. Cells that are already in the normalised form (and any
# single-entry / blank cells) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$alreadyNormalised = "backup@backdoor.com, System"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Text

    if ($val -ne $null -and $val -ne "" -and $val -ne $alreadyNormalised) {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $newVal = $reversed -join ", "
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
